# Update price list: move the date forward one month and refresh unit
# prices in the PITONES ROSCA W. table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Date in A1 moves from 24-Apr-2024 to 24-May-2024 (serial 45406 -> 45436)
$ws.Range("A1").Value = Get-Date -Year 2024 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Refreshed prices in column D for rows 33-38
$ws.Range("D33").Value = 4282.249
$ws.Range("D34").Value = 4282.249
$ws.Range("D35").Value = 5601.891
$ws.Range("D36").Value = 5601.932
$ws.Range("D37").Value = 5601.891
$ws.Range("D38").Value = 7061.332
